$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "numeric-looking" string (must be forced to stay text,
# matching the workbooks original inlineStr/text cell type for the Price column).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$ws.Range("D2").Value = "26.340.90"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "1.712.90"
$ws.Range("E3").Value = "  -1.39%  "

Set-TextValue $ws.Range("D4") "0.9966"
$ws.Range("E4").Value = "  -0.33%  "

Set-TextValue $ws.Range("D5") "240.06"
$ws.Range("E5").Value = "  -2.69%  "

Set-TextValue $ws.Range("D6") "0.9974"
$ws.Range("E6").Value = "  -0.28%  "

Set-TextValue $ws.Range("D7") "0.4855"
$ws.Range("E7").Value = "  -1.19%  "

$ws.Range("E8").Value = "  -3.14%  "

Set-TextValue $ws.Range("D9") "0.06170"
$ws.Range("E9").Value = "  -2.12%  "

$ws.Range("D10").Value = "1.714.82"
$ws.Range("E10").Value = "  -1.27%  "

Set-TextValue $ws.Range("D11") "0.06941"
$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("E13").Value = "  -3.01%  "

Set-TextValue $ws.Range("D14") "0.5960"
$ws.Range("E14").Value = "  -2.50%  "

Set-TextValue $ws.Range("D15") "76.35"
$ws.Range("E15").Value = "  -1.42%  "

Set-TextValue $ws.Range("D16") "0.9973"
$ws.Range("E16").Value = "  -0.27%  "

Set-TextValue $ws.Range("D17") "0.9968"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").Value = "26.244.93"
$ws.Range("E18").Value = "  -1.07%  "

Set-TextValue $ws.Range("D19") "0.000007096"
$ws.Range("E19").Value = "  -4.06%  "

$ws.Range("E20").Value = "  -2.84%  "

$ws.Range("D21").Value = "1.932.01"
$ws.Range("E21").Value = "  -1.29%  "

Set-TextValue $ws.Range("D22") "4.397"
$ws.Range("E22").Value = "  -4.20%  "

Set-TextValue $ws.Range("D23") "8.408"
$ws.Range("E23").Value = "  -3.52%  "

Set-TextValue $ws.Range("D24") "5.033"
$ws.Range("E24").Value = "  -4.22%  "

Set-TextValue $ws.Range("D25") "135.98"
$ws.Range("E25").Value = "  -2.88%  "

Set-TextValue $ws.Range("D26") "15.14"
$ws.Range("E26").Value = "  -2.08%  "

Set-TextValue $ws.Range("D27") "1.394"
$ws.Range("E27").Value = "  -1.81%  "

Set-TextValue $ws.Range("D28") "1.725"
$ws.Range("E28").Value = "  -2.20%  "

Set-TextValue $ws.Range("D29") "105.34"
$ws.Range("E29").Value = "  -2.43%  "

Set-TextValue $ws.Range("D30") "3.861"
$ws.Range("E30").Value = "  -4.51%  "

Set-TextValue $ws.Range("D31") "0.07938"
$ws.Range("E31").Value = "  -1.53%  "

Set-TextValue $ws.Range("D32") "3.600"
$ws.Range("E32").Value = "  -3.12%  "

Set-TextValue $ws.Range("D33") "0.04420"
$ws.Range("E33").Value = "  -3.74%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D34") "2.599"
$ws.Range("E34").Value = "  -0.48%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "0.9877"
$ws.Range("E35").Value = "  -2.09%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "0.6170"
$ws.Range("E36").Value = "  -3.19%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D37") "0.9318"
$ws.Range("E37").Value = "  +3.97%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D38") "1.979"
$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D39") "2.367"
$ws.Range("E39").Value = "  -1.43%  "

$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D40") "0.9963"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.01469"
$ws.Range("E41").Value = "  -2.60%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D42") "99.84"
$ws.Range("E42").Value = "  -2.45%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "5.373"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D44") "0.3798"
$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D45") "6.834"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D46") "0.1146"
$ws.Range("E46").Value = "  -3.45%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.05340"
$ws.Range("E47").Value = "  -1.10%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D48") "30.53"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "7.699"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D50") "50.99"
$ws.Range("E50").Value = "  -1.58%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.209"
$ws.Range("E51").Value = "  -4.82%  "

